# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Column letters -> indices: A=1, B=2, C=3, D=4, E=5
# Price values in column D are stored as *text* (e.g. "212.48", "29.679.72"
# with thousands separators as literal dots). Excel's COM layer auto-detects
# plain numeric-looking strings as numbers, so a leading apostrophe is used
# to force text entry, exactly as typing '212.48 into a cell would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

Set-TextValue 2 4 "29.679.72"
$ws.Cells.Item(2, 5).Value = "  +0.72%  "

Set-TextValue 3 4 "1.608.13"
$ws.Cells.Item(3, 5).Value = "  +0.33%  "

$ws.Cells.Item(4, 5).Value = "  -0.24%  "

Set-TextValue 5 4 "212.48"
$ws.Cells.Item(5, 5).Value = "  -0.22%  "

$ws.Cells.Item(6, 5).Value = "  +0.03%  "

Set-TextValue 7 4 "0.995"
$ws.Cells.Item(7, 5).Value = "  -0.24%  "

Set-TextValue 8 4 "28.85"
$ws.Cells.Item(8, 5).Value = "  +7.28%  "

$ws.Cells.Item(9, 5).Value = "  +3.60%  "

$ws.Cells.Item(10, 5).Value = "  +1.50%  "

Set-TextValue 11 4 "0.0905"
$ws.Cells.Item(11, 5).Value = "  -0.46%  "

Set-TextValue 12 4 "1.839.30"
$ws.Cells.Item(12, 5).Value = "  +0.45%  "

Set-TextValue 13 4 "1.604.95"
$ws.Cells.Item(13, 5).Value = "  -0.70%  "

Set-TextValue 14 4 "0.565"
$ws.Cells.Item(14, 5).Value = "  +5.57%  "

Set-TextValue 15 4 "3.86"
$ws.Cells.Item(15, 5).Value = "  +3.34%  "

Set-TextValue 16 4 "29.694.70"
$ws.Cells.Item(16, 5).Value = "  +0.57%  "

Set-TextValue 17 4 "8.70"
$ws.Cells.Item(17, 5).Value = "  +14.39%  "

Set-TextValue 18 4 "64.66"
$ws.Cells.Item(18, 5).Value = "  +1.89%  "

Set-TextValue 19 4 "241.63"
$ws.Cells.Item(19, 5).Value = "  +0.78%  "

Set-TextValue 20 4 "0.0₃0704"
$ws.Cells.Item(20, 5).Value = "  +1.73%  "

$ws.Cells.Item(21, 5).Value = "  -0.20%  "

# Row 22 / 23 swap content (Uniswap <-> Avalanche)
$ws.Cells.Item(22, 2).Value = "Avalanche"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue 22 4 "9.68"
$ws.Cells.Item(22, 5).Value = "  +5.68%  "

$ws.Cells.Item(23, 2).Value = "Uniswap"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue 23 4 "4.07"
$ws.Cells.Item(23, 5).Value = "  +1.41%  "

Set-TextValue 24 4 "2.12"
$ws.Cells.Item(24, 5).Value = "  +2.74%  "

Set-TextValue 25 4 "156.79"
$ws.Cells.Item(25, 5).Value = "  +1.56%  "

Set-TextValue 26 4 "15.59"
$ws.Cells.Item(26, 5).Value = "  +2.05%  "

Set-TextValue 27 4 "0.110"
$ws.Cells.Item(27, 5).Value = "  +1.17%  "

Set-TextValue 28 4 "6.57"
$ws.Cells.Item(28, 5).Value = "  +3.32%  "

$ws.Cells.Item(30, 5).Value = "  +2.12%  "

# Row 31 / 32 swap content (PancakeSwap <-> Filecoin)
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 31 4 "3.28"
$ws.Cells.Item(31, 5).Value = "  +1.16%  "

$ws.Cells.Item(32, 2).Value = "PancakeSwap"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue 32 4 "1.07"
$ws.Cells.Item(32, 5).Value = "  +0.33%  "

$ws.Cells.Item(33, 5).Value = "  +2.61%  "

Set-TextValue 34 4 "1.430.00"
$ws.Cells.Item(34, 5).Value = "  +0.42%  "

Set-TextValue 35 4 "1.61"
$ws.Cells.Item(35, 5).Value = "  +6.54%  "

$ws.Cells.Item(36, 5).Value = "  +0.93%  "

Set-TextValue 37 4 "2.87"
$ws.Cells.Item(37, 5).Value = "  +2.16%  "

$ws.Cells.Item(38, 5).Value = "  -0.62%  "

$ws.Cells.Item(39, 5).Value = "  +3.24%  "

Set-TextValue 40 4 "0.554"
$ws.Cells.Item(40, 5).Value = "  +3.79%  "

Set-TextValue 41 4 "0.0498"
$ws.Cells.Item(41, 5).Value = "  +5.76%  "

# Row 42 / 43 swap content (RenderToken <-> ARBITRUM)
$ws.Cells.Item(42, 2).Value = "ARBITRUM"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue 42 4 "0.826"
$ws.Cells.Item(42, 5).Value = "  +4.24%  "

$ws.Cells.Item(43, 2).Value = "RenderToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 43 4 "1.97"
$ws.Cells.Item(43, 5).Value = "  +0.93%  "

Set-TextValue 44 4 "54.35"
$ws.Cells.Item(44, 5).Value = "  +2.37%  "

Set-TextValue 45 4 "68.96"
$ws.Cells.Item(45, 5).Value = "  +5.63%  "

$ws.Cells.Item(46, 5).Value = "  -0.20%  "

Set-TextValue 47 4 "1.00"
$ws.Cells.Item(47, 5).Value = "  +19.72%  "

Set-TextValue 48 4 "5.45"
$ws.Cells.Item(48, 5).Value = "  +3.15%  "

Set-TextValue 49 4 "1.747.75"
$ws.Cells.Item(49, 5).Value = "  +0.23%  "

Set-TextValue 50 4 "87.02"
$ws.Cells.Item(50, 5).Value = "  +0.67%  "

$ws.Cells.Item(51, 5).Value = "  -1.15%  "
